$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh "last updated" timestamp (19:39 -> 21:04)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 21:04"

# Updated COVID-19 country data. The source sheet is kept sorted by
# "Casos totales" (column B) descending, so a handful of countries also
# change row/rank along with their refreshed figures.
$data = @(
    @(4, $null, 1308977, 16354, 220983, 1010080, 16798, 986, 77914),
    @(9, $null, 174791, 0, 55782, 92779, 2868, 243, 26230),
    @(10, $null, 170090, 660, 141700, 20951, 1712, 47, 7439),
    @(11, $null, 141088, 5395, 55350, 76101, 8318, 449, 9637),
    @(15, $null, 66313, 1391, 29942, 31804, 502, 159, 4567),
    @(16, $null, 59693, 3342, 17883, 39825, 0, 96, 1985),
    @(35, $null, 15575, 98, 5146, 9839, 300, 13, 590),
    @(57, $null, 5371, 0, 1659, 3427, 148, 3, 285),
    @(92, "Sudan", 1111, 181, 102, 950, 0, 7, 59),
    @(93, "Hong Kong", 1045, 0, 960, 81, 1, 0, 4),
    @(94, "Tunez", 1026, 0, 600, 382, 22, 0, 44),
    @(108, "Maldivas", 744, 96, 20, 721, 2, 0, 3),
    @(109, "El Salvador", 742, 47, 257, 469, 4, 1, 16),
    @(110, "Burkina Faso", 736, 0, 562, 126, 0, 0, 48),
    @(158, $null, 105, 4, 30, 75, 0, 0, 0),
    @(164, $null, 92, 0, 31, 50, 1, 0, 11),
    @(165, "Sudan del Sur", 90, 16, 2, 88, 0, 0, 0),
    @(166, "Barbados", 83, 1, 53, 23, 4, 0, 7),
    @(167, "Mozambique", 82, 1, 27, 55, 0, 0, 0),
    @(168, "Liechtenstein", 82, 0, 55, 26, 0, 0, 1),
    @(169, "Islas Caimanes", 80, 0, 35, 44, 3, 0, 1),
    @(170, "San Martin (Parte Holandesa)", 76, 0, 44, 18, 7, 0, 14),
    @(178, $null, 39, 0, 37, 2, 0, 0, 0),
    @(192, "Belice", 18, 0, 16, 0, 0, 0, 2),
    @(193, "Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $countryName = $row[1]
    if ($countryName -ne $null) {
        $ws.Cells.Item($r, 1).Value = $countryName
    }
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}